# Burndown Sprint 20 - update task list and daily progress numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update task names (US103 -> dropped, US104 keeps its text but moves to a
# different shared-string slot, and a brand new US116 task replaces the old
# US103 row's underlying string) ---
$ws.Range("A3").Value = "#US104 Agregar boton de terminar al examen"
$ws.Range("A7").Value = "#US116 Ejecutar ciclos de test"

# --- Row 3: US104 ---
$ws.Range("B3").Value = 5
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 5
$ws.Range("W3").Value = 0

# --- Row 4: US105 ---
$ws.Range("B4").Value = 5
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = 0

# --- Row 5: US107 ---
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 0
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0

# --- Row 6: US112 ---
$ws.Range("B6").Value = 100
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 7
$ws.Range("I6").Value = 8
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 12
$ws.Range("P6").Value = 7
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 12
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 6

# --- Row 7: US116 ---
$ws.Range("B7").Value = 50
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 0
$ws.Range("L7").Value = 10
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 5
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 5
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0

# Move the active selection on the frozen pane to I4, matching the saved view
$ws.Range("I4").Select()

$excel.Calculate()
$wb.Save()
